$d = $word.ActiveDocument

$replacements = @(
    @("452÷4=", "918÷4="),
    @("627÷2=", "681÷2="),
    @("433÷2=", "148÷5="),
    @("236÷5=", "682÷3="),
    @("458÷3=", "430÷6="),
    @("573÷9=", "314÷7="),
    @("797÷2=", "335÷7="),
    @("586÷6=", "804÷2="),
    @("707÷8=", "935÷8="),
    @("478÷6=", "309÷7="),
    @("210÷2=", "328÷9="),
    @("429÷3=", "532÷3="),
    @("232÷5=", "364÷4="),
    @("201÷9=", "688÷4="),
    @("795÷4=", "436÷9="),
    @("280÷7=", "134÷5="),
    @("961÷6=", "595÷6="),
    @("723÷6=", "379÷4="),
    @("795÷7=", "149÷9="),
    @("332÷9=", "482÷6="),
    @("771÷7=", "549÷5="),
    @("152÷3=", "158÷4="),
    @("789÷6=", "333÷9="),
    @("279÷4=", "451÷7="),
    @("625÷7=", "380÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
